$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22-49 down to 23-50.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = "Macroferia Regional de Talca"
$ws.Range("C22").Value = "Maule"
$ws.Range("D22").Value = 44495
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = 100112022
$ws.Range("G22").Value = "Arveja Verde"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 17000
$ws.Range("M22").Value = 17000
$ws.Range("N22").Value = '$/saco 25 kilos'
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 680
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
